$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/percentage/coin-name/link updates - safe to assign directly as strings
# since none of these look like a bare number to Excel's type inference.
$plainUpdates = @{
    'D2' = '36.320.69'
    'E2' = '  -1.68%  '
    'D3' = '2.047.88'
    'E3' = '  -3.10%  '
    'E4' = '  +0.27%  '
    'E5' = '  -0.85%  '
    'E6' = '  +0.83%  '
    'E7' = '  +0.10%  '
    'E8' = '  -0.06%  '
    'E9' = '  +6.23%  '
    'E10' = '  -1.14%  '
    'E11' = '  -3.97%  '
    'E12' = '  -3.73%  '
    'E13' = '  +1.80%  '
    'E14' = '  -6.22%  '
    'D15' = '2.345.70'
    'E15' = '  -3.07%  '
    'E16' = '  -3.83%  '
    'D17' = '2.064.62'
    'E17' = '  -4.50%  '
    'E18' = '  +1.85%  '
    'D19' = '36.286.57'
    'E19' = '  -1.71%  '
    'E20' = '  -3.14%  '
    'D21' = '0.0₃0855'
    'E21' = '  -3.84%  '
    'E22' = '  -1.06%  '
    'E23' = '  -7.16%  '
    'E24' = '  +0.20%  '
    'E25' = '  -3.35%  '
    'E27' = '  -6.78%  '
    'E28' = '  -2.33%  '
    'E29' = '  -5.03%  '
    'E30' = '  -2.57%  '
    'E31' = '  -1.97%  '
    'E32' = '  -8.27%  '
    'E33' = '  -2.90%  '
    'E34' = '  -7.01%  '
    'B35' = 'Kaspa'
    'C35' = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    'E35' = '  +3.01%  '
    'B36' = 'BinanceUSD'
    'C36' = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
    'E36' = '  +0.11%  '
    'E37' = '  -0.83%  '
    'E38' = '  -9.38%  '
    'E39' = '  +1.62%  '
    'E40' = '  -6.90%  '
    'B41' = 'VeChain'
    'C41' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'E41' = '  -3.45%  '
    'B42' = 'HuobiToken'
    'C42' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'E42' = '  -0.05%  '
    'E43' = '  -7.41%  '
    'E44' = '  -5.01%  '
    'E45' = '  -5.29%  '
    'D46' = '1.376.30'
    'E46' = '  +0.74%  '
    'E47' = '  -2.30%  '
    'E48' = '  +4.46%  '
    'E50' = '  -9.29%  '
    'B51' = 'MultiversX'
    'C51' = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
    'E51' = '  -0.51%  '
}
foreach ($cellRef in $plainUpdates.Keys) {
    $ws.Range($cellRef).Value = $plainUpdates[$cellRef]
}

# Price updates in column D whose new text parses as a plain number
# (e.g. "56.55"). A direct .Value assignment would be auto-coerced to a
# numeric cell, changing the stored type from the original text cell.
# Force the cell to Text first, write the value, then restore the
# cell's formatting to the sheet default so no stray number-format/
# quote-prefix style lingers on the cell - matching the original
# plain (unstyled) text cells.
$numericLookingUpdates = @{
    'D5' = '244.96'
    'D8' = '56.55'
    'D9' = '63.80'
    'D10' = '0.368'
    'D11' = '0.0746'
    'D13' = '0.912'
    'D14' = '14.33'
    'D16' = '5.39'
    'D18' = '17.79'
    'D20' = '71.40'
    'D22' = '236.79'
    'D23' = '5.20'
    'D25' = '2.35'
    'D27' = '9.30'
    'D28' = '164.75'
    'D29' = '19.96'
    'D32' = '4.96'
    'D33' = '0.0598'
    'D34' = '4.41'
    'D35' = '0.0874'
    'D36' = '1.00'
    'D37' = '1.83'
    'D38' = '2.22'
    'D39' = '5.04'
    'D41' = '0.0215'
    'D42' = '2.87'
    'D44' = '93.47'
    'D45' = '0.0905'
    'D47' = '15.86'
    'D48' = '7.42'
    'D49' = '2.95'
    'D50' = '2.28'
    'D51' = '45.74'
}
foreach ($cellRef in $numericLookingUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingUpdates[$cellRef]
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}
